$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row for Pomelo "Start Ruby" / "Primera" needs to be
# recorded. Insert a fresh row at 568 (pushing the existing rows 568..651 down
# to 569..652, same as Excel's normal "insert row above" behaviour) and fill
# it in with the new week's data. All of the descriptive/static columns match
# the surrounding rows for this market+product+variety+quality combination.

$ws.Rows.Item(568).Insert()

$ws.Range("A568").Value = 4
$ws.Range("B568").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C568").Value = "Los Lagos"
$ws.Range("D568").Value = 45131
$ws.Range("E568").Value = 10
$ws.Range("F568").Value = "Fruta"
$ws.Range("G568").Value = 100102
$ws.Range("H568").Value = "Cítricos"
$ws.Range("I568").Value = 100102006
$ws.Range("J568").Value = "Pomelo"
$ws.Range("K568").Value = "Start Ruby"
$ws.Range("L568").Value = "Primera"
$ws.Range("M568").Value = 100
$ws.Range("N568").Value = 14000
$ws.Range("O568").Value = 14000
$ws.Range("P568").Value = 14000
$ws.Range("Q568").Value = "$/caja 14 kilos empedrada"
$ws.Range("R568").Value = "Región de O'Higgins"
$ws.Range("S568").Value = 1000
$ws.Range("T568").Value = 14
